$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DVD (math) rows to append below the existing data (rows 52-57).
$newItems = @(
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 1",
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 2",
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 3",
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 4",
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 5",
    "วีดิทัศน์ (DVD) ประกอบการเรียนรู้วิชาคณิตศาสตร์  ชั้นประถมศึกษาปีที่ 6"
)

$startRow = 52
for ($i = 0; $i -lt $newItems.Count; $i++) {
    $row = $startRow + $i

    # Copy formatting (style) from the existing "math" category band (A30:B30)
    # down onto the new row before writing the values.
    $ws.Range("A30:B30").Copy() | Out-Null
    $ws.Range("A" + $row + ":B" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $row).Value = $newItems[$i]
    $ws.Range("B" + $row).Value = ($startRow - 1 + $i)
}

$excel.CutCopyMode = $false

# Match the final view state recorded in the saved workbook.
$excel.Goto($ws.Range("A34"), $true) | Out-Null
$ws.Range("H48").Select() | Out-Null
